$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10 (Id 45000007): "first clear" -> "lost" entry
$ws.Range("B10").Value = "迷失"
$ws.Range("C10").Value = "连续3次迷失在森林中，没能进入更深处"
$ws.Range("D10").Value = 4
$ws.Range("G10").Value = "innerforest"
$ws.Range("J10").Value = 3
$ws.Range("L10").Value = "box"

# Row 11 (Id 45000008): "first clear" -> "lost track" entry
$ws.Range("B11").Value = "走失"
$ws.Range("C11").Value = "没有找到穷奇"
$ws.Range("D11").Value = 3
$ws.Range("G11").Value = "bossqiongqi2"
$ws.Range("I11").Value = 1
$ws.Range("L11").Value = "duelist2"

# Update selection to L11 to match saved cursor position
$ws.Range("L11").Select() | Out-Null
